# The author opened the workbook, clicked into Sheet1!C3 (which held "A"),
# briefly typed a replacement value and then reverted it before saving -
# net effect on visible content is nil, but the active selection ends up
# on C3 instead of the original G9.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Touch the cell (mirrors the transient "A B    " edit in the source commit)
# and restore its original value so the saved data is unchanged.
$original = $ws.Range("C3").Value2
$ws.Range("C3").Value = "A B    "
$ws.Range("C3").Value = $original

# Leave the selection on C3, matching the saved sheetView.
$ws.Range("C3").Select()
